$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 15.953202
$ws.Range("H2").Value = 47.859606
$ws.Range("I2").Value = 0.6210379196599995
$ws.Range("J2").Value = 0.6210379196599995
$ws.Range("M2").Value = 2.231113333333334
$ws.Range("N2").Value = 6.69334
$ws.Range("O2").Value = 0.01598125358798882
$ws.Range("P2").Value = 0.01598125358798882
$ws.Range("Q2").Value = 35.59340169156
$ws.Range("R2").Value = 320.34061522404
$ws.Range("S2").Value = 0.00992496448184348
$ws.Range("T2").Value = 0.00992496448184348

# Row 3
$ws.Range("G3").Value = 15.953202
$ws.Range("H3").Value = 47.859606
$ws.Range("I3").Value = 0.6210379196599995
$ws.Range("J3").Value = 0.6210379196599995
$ws.Range("O3").Value = 0.1634493267640196
$ws.Range("P3").Value = 0.1634493267640195
$ws.Range("Q3").Value = 364.033867036518
$ws.Range("R3").Value = 3276.304803328662
$ws.Range("S3").Value = 0.1015082298633542
$ws.Range("T3").Value = 0.1015082298633542

# Row 4
$ws.Range("G4").Value = 15.953202
$ws.Range("H4").Value = 47.859606
$ws.Range("I4").Value = 0.6210379196599995
$ws.Range("J4").Value = 0.6210379196599995
$ws.Range("M4").Value = 58.02175166666666
$ws.Range("N4").Value = 174.065255
$ws.Range("O4").Value = 0.4156043142904646
$ws.Range("P4").Value = 0.4156043142904646
$ws.Range("Q4").Value = 925.6327247321699
$ws.Range("R4").Value = 8330.694522589529
$ws.Range("S4").Value = 0.2581060387486707
$ws.Range("T4").Value = 0.2581060387486707

# Row 5
$ws.Range("G5").Value = 15.953202
$ws.Range("H5").Value = 47.859606
$ws.Range("I5").Value = 0.6210379196599995
$ws.Range("J5").Value = 0.6210379196599995
$ws.Range("M5").Value = 15.16934033333333
$ws.Range("N5").Value = 45.508021
$ws.Range("O5").Value = 0.1086565487318021
$ws.Range("P5").Value = 0.1086565487318021
$ws.Range("Q5").Value = 241.999550544414
$ws.Range("R5").Value = 2177.995954899726
$ws.Range("S5").Value = 0.06747983698183375
$ws.Range("T5").Value = 0.06747983698183375

# Row 6
$ws.Range("G6").Value = 15.953202
$ws.Range("H6").Value = 47.859606
$ws.Range("I6").Value = 0.6210379196599995
$ws.Range("J6").Value = 0.6210379196599995
$ws.Range("M6").Value = 41.36709099999999
$ws.Range("N6").Value = 124.101273
$ws.Range("O6").Value = 0.2963085566257249
$ws.Range("P6").Value = 0.2963085566257249
$ws.Range("Q6").Value = 659.9375588753819
$ws.Range("R6").Value = 5939.438029878437
$ws.Range("S6").Value = 0.1840188495842974
$ws.Range("T6").Value = 0.1840188495842974

# Row 7
$ws.Range("G7").Value = 0.7397413333333333
$ws.Range("I7").Value = 0.02879719185777549
$ws.Range("J7").Value = 0.02879719185777549
$ws.Range("M7").Value = 2.231113333333334
$ws.Range("N7").Value = 6.69334
$ws.Range("O7").Value = 0.01598125358798882
$ws.Range("P7").Value = 0.01598125358798882
$ws.Range("Q7").Value = 1.650446752017778
$ws.Range("R7").Value = 14.85402076816
$ws.Range("S7").Value = 0.0004602152257010768
$ws.Range("T7").Value = 0.0004602152257010768

# Row 8
$ws.Range("G8").Value = 0.7397413333333333
$ws.Range("I8").Value = 0.02879719185777549
$ws.Range("J8").Value = 0.02879719185777549
$ws.Range("O8").Value = 0.1634493267640196
$ws.Range("P8").Value = 0.1634493267640195
$ws.Range("S8").Value = 0.004706881621847709
$ws.Range("T8").Value = 0.004706881621847708

# Row 9
$ws.Range("G9").Value = 0.7397413333333333
$ws.Range("I9").Value = 0.02879719185777549
$ws.Range("J9").Value = 0.02879719185777549
$ws.Range("M9").Value = 58.02175166666666
$ws.Range("N9").Value = 174.065255
$ws.Range("O9").Value = 0.4156043142904646
$ws.Range("P9").Value = 0.4156043142904646
$ws.Range("Q9").Value = 42.92108794023554
$ws.Range("R9").Value = 386.2897914621199
$ws.Range("S9").Value = 0.01196823717554173
$ws.Range("T9").Value = 0.01196823717554173

# Row 10
$ws.Range("G10").Value = 0.7397413333333333
$ws.Range("I10").Value = 0.02879719185777549
$ws.Range("J10").Value = 0.02879719185777549
$ws.Range("M10").Value = 15.16934033333333
$ws.Range("N10").Value = 45.508021
$ws.Range("O10").Value = 0.1086565487318021
$ws.Range("P10").Value = 0.1086565487318021
$ws.Range("Q10").Value = 11.22138804396711
$ws.Range("R10").Value = 100.992492395704
$ws.Range("S10").Value = 0.003129003480433437
$ws.Range("T10").Value = 0.003129003480433437

# Row 11
$ws.Range("G11").Value = 0.7397413333333333
$ws.Range("I11").Value = 0.02879719185777549
$ws.Range("J11").Value = 0.02879719185777549
$ws.Range("M11").Value = 41.36709099999999
$ws.Range("N11").Value = 124.101273
$ws.Range("O11").Value = 0.2963085566257249
$ws.Range("P11").Value = 0.2963085566257249
$ws.Range("Q11").Value = 30.60094705246133
$ws.Range("R11").Value = 275.408523472152
$ws.Range("S11").Value = 0.008532854354251532
$ws.Range("T11").Value = 0.008532854354251532

# Row 12
$ws.Range("G12").Value = 5.607355000000001
$ws.Range("H12").Value = 16.822065
$ws.Range("I12").Value = 0.2182872180766656
$ws.Range("J12").Value = 0.2182872180766656
$ws.Range("M12").Value = 2.231113333333334
$ws.Range("N12").Value = 6.69334
$ws.Range("O12").Value = 0.01598125358798882
$ws.Range("P12").Value = 0.01598125358798882
$ws.Range("Q12").Value = 12.51064450523334
$ws.Range("R12").Value = 112.5958005471
$ws.Range("S12").Value = 0.00348850338709981
$ws.Range("T12").Value = 0.003488503387099809

# Row 13
$ws.Range("G13").Value = 5.607355000000001
$ws.Range("H13").Value = 16.822065
$ws.Range("I13").Value = 0.2182872180766656
$ws.Range("J13").Value = 0.2182872180766656
$ws.Range("O13").Value = 0.1634493267640196
$ws.Range("P13").Value = 0.1634493267640195
$ws.Range("Q13").Value = 127.953443107945
$ws.Range("R13").Value = 1151.580987971505
$ws.Range("S13").Value = 0.03567889883582172
$ws.Range("T13").Value = 0.0356788988358217

# Row 14
$ws.Range("G14").Value = 5.607355000000001
$ws.Range("H14").Value = 16.822065
$ws.Range("I14").Value = 0.2182872180766656
$ws.Range("J14").Value = 0.2182872180766656
$ws.Range("M14").Value = 58.02175166666666
$ws.Range("N14").Value = 174.065255
$ws.Range("O14").Value = 0.4156043142904646
$ws.Range("P14").Value = 0.4156043142904646
$ws.Range("Q14").Value = 325.3485593168417
$ws.Range("R14").Value = 2928.137033851575
$ws.Range("S14").Value = 0.0907211095871257
$ws.Range("T14").Value = 0.09072110958712569

# Row 15
$ws.Range("G15").Value = 5.607355000000001
$ws.Range("H15").Value = 16.822065
$ws.Range("I15").Value = 0.2182872180766656
$ws.Range("J15").Value = 0.2182872180766656
$ws.Range("M15").Value = 15.16934033333333
$ws.Range("N15").Value = 45.508021
$ws.Range("O15").Value = 0.1086565487318021
$ws.Range("P15").Value = 0.1086565487318021
$ws.Range("Q15").Value = 85.05987636481835
$ws.Range("R15").Value = 765.5388872833651
$ws.Range("S15").Value = 0.02371833574847673
$ws.Range("T15").Value = 0.02371833574847673

# Row 16
$ws.Range("G16").Value = 5.607355000000001
$ws.Range("H16").Value = 16.822065
$ws.Range("I16").Value = 0.2182872180766656
$ws.Range("J16").Value = 0.2182872180766656
$ws.Range("M16").Value = 41.36709099999999
$ws.Range("N16").Value = 124.101273
$ws.Range("O16").Value = 0.2963085566257249
$ws.Range("P16").Value = 0.2963085566257249
$ws.Range("Q16").Value = 231.959964554305
$ws.Range("R16").Value = 2087.639680988745
$ws.Range("S16").Value = 0.06468037051814163
$ws.Range("T16").Value = 0.06468037051814163

# Row 17
$ws.Range("G17").Value = 1.453021
$ws.Range("H17").Value = 4.359063
$ws.Range("I17").Value = 0.05656426459479998
$ws.Range("J17").Value = 0.05656426459479998
$ws.Range("M17").Value = 2.231113333333334
$ws.Range("N17").Value = 6.69334
$ws.Range("O17").Value = 0.01598125358798882
$ws.Range("P17").Value = 0.01598125358798882
$ws.Range("Q17").Value = 3.241854526713333
$ws.Range("R17").Value = 29.17669074042
$ws.Range("S17").Value = 0.000903967856507596
$ws.Range("T17").Value = 0.0009039678565075961

# Row 18
$ws.Range("G18").Value = 1.453021
$ws.Range("H18").Value = 4.359063
$ws.Range("I18").Value = 0.05656426459479998
$ws.Range("J18").Value = 0.05656426459479998
$ws.Range("O18").Value = 0.1634493267640196
$ws.Range("P18").Value = 0.1634493267640195
$ws.Range("Q18").Value = 33.156281323039
$ws.Range("R18").Value = 298.406531907351
$ws.Range("S18").Value = 0.009245390966921925
$ws.Range("T18").Value = 0.009245390966921924

# Row 19
$ws.Range("G19").Value = 1.453021
$ws.Range("H19").Value = 4.359063
$ws.Range("I19").Value = 0.05656426459479998
$ws.Range("J19").Value = 0.05656426459479998
$ws.Range("M19").Value = 58.02175166666666
$ws.Range("N19").Value = 174.065255
$ws.Range("O19").Value = 0.4156043142904646
$ws.Range("P19").Value = 0.4156043142904646
$ws.Range("Q19").Value = 84.30682362845165
$ws.Range("R19").Value = 758.7614126560649
$ws.Range("S19").Value = 0.02350835240026625
$ws.Range("T19").Value = 0.02350835240026625

# Row 20
$ws.Range("G20").Value = 1.453021
$ws.Range("H20").Value = 4.359063
$ws.Range("I20").Value = 0.05656426459479998
$ws.Range("J20").Value = 0.05656426459479998
$ws.Range("M20").Value = 15.16934033333333
$ws.Range("N20").Value = 45.508021
$ws.Range("O20").Value = 0.1086565487318021
$ws.Range("P20").Value = 0.1086565487318021
$ws.Range("Q20").Value = 22.04137006048033
$ws.Range("R20").Value = 198.372330544323
$ws.Range("S20").Value = 0.006146077772423433
$ws.Range("T20").Value = 0.006146077772423434

# Row 21
$ws.Range("G21").Value = 1.453021
$ws.Range("H21").Value = 4.359063
$ws.Range("I21").Value = 0.05656426459479998
$ws.Range("J21").Value = 0.05656426459479998
$ws.Range("M21").Value = 41.36709099999999
$ws.Range("N21").Value = 124.101273
$ws.Range("O21").Value = 0.2963085566257249
$ws.Range("P21").Value = 0.2963085566257249
$ws.Range("Q21").Value = 60.10725193191099
$ws.Range("R21").Value = 540.965267387199
$ws.Range("S21").Value = 0.01676047559868078
$ws.Range("T21").Value = 0.01676047559868078

# Row 22
$ws.Range("G22").Value = 1.934648333333333
$ws.Range("H22").Value = 5.803945
$ws.Range("I22").Value = 0.07531340581075942
$ws.Range("J22").Value = 0.07531340581075942
$ws.Range("M22").Value = 2.231113333333334
$ws.Range("N22").Value = 6.69334
$ws.Range("O22").Value = 0.01598125358798882
$ws.Range("P22").Value = 0.01598125358798882
$ws.Range("Q22").Value = 4.316419691811111
$ws.Range("R22").Value = 38.8477772263
$ws.Range("S22").Value = 0.001203602636836857
$ws.Range("T22").Value = 0.001203602636836857

# Row 23
$ws.Range("G23").Value = 1.934648333333333
$ws.Range("H23").Value = 5.803945
$ws.Range("I23").Value = 0.07531340581075942
$ws.Range("J23").Value = 0.07531340581075942
$ws.Range("O23").Value = 0.1634493267640196
$ws.Range("P23").Value = 0.1634493267640195
$ws.Range("Q23").Value = 44.14646753291833
$ws.Range("R23").Value = 397.318207796265
$ws.Range("S23").Value = 0.01230992547607403
$ws.Range("T23").Value = 0.01230992547607403

# Row 24
$ws.Range("G24").Value = 1.934648333333333
$ws.Range("H24").Value = 5.803945
$ws.Range("I24").Value = 0.07531340581075942
$ws.Range("J24").Value = 0.07531340581075942
$ws.Range("M24").Value = 58.02175166666666
$ws.Range("N24").Value = 174.065255
$ws.Range("O24").Value = 0.4156043142904646
$ws.Range("P24").Value = 0.4156043142904646
$ws.Range("Q24").Value = 112.2516851589972
$ws.Range("R24").Value = 1010.265166430975
$ws.Range("S24").Value = 0.03130057637886016
$ws.Range("T24").Value = 0.03130057637886016

# Row 25
$ws.Range("G25").Value = 1.934648333333333
$ws.Range("H25").Value = 5.803945
$ws.Range("I25").Value = 0.07531340581075942
$ws.Range("J25").Value = 0.07531340581075942
$ws.Range("M25").Value = 15.16934033333333
$ws.Range("N25").Value = 45.508021
$ws.Range("O25").Value = 0.1086565487318021
$ws.Range("P25").Value = 0.1086565487318021
$ws.Range("Q25").Value = 29.34733899364944
$ws.Range("R25").Value = 264.126050942845
$ws.Range("S25").Value = 0.008183294748634769
$ws.Range("T25").Value = 0.008183294748634769

# Row 26
$ws.Range("G26").Value = 1.934648333333333
$ws.Range("H26").Value = 5.803945
$ws.Range("I26").Value = 0.07531340581075942
$ws.Range("J26").Value = 0.07531340581075942
$ws.Range("M26").Value = 41.36709099999999
$ws.Range("N26").Value = 124.101273
$ws.Range("O26").Value = 0.2963085566257249
$ws.Range("P26").Value = 0.2963085566257249
$ws.Range("Q26").Value = 80.03077365799832
$ws.Range("R26").Value = 720.2769629219849
$ws.Range("S26").Value = 0.02231600657035361
$ws.Range("T26").Value = 0.02231600657035361

Write-Output "Applied 283 cell updates"
